$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- View tweaks: zoom 130% -> 120%, selection moves to C20 ---
$excel.ActiveWindow.Zoom = 120
$ws.Range("C20").Select()

# --- Row 6 / 7 : BMsolveVFI <-> BMsolveGSSA ---
# Row 6 becomes BMsolveGSSA and now carries the new result (F6 / F6/F2)
$ws.Range("A6").Value = "BMsolveGSSA"
$ws.Range("F6").Value = 3617.7648222972498
$ws.Range("B6").Formula = "=F6/F2"
# Row 7 becomes BMsolveVFI, no numbers
$ws.Range("A7").Value = "BMsolveVFI"

# --- Row 10 / 11 : BMsimVFI <-> BMsimGSSA ---
$ws.Range("A10").Value = "BMsimGSSA"
$ws.Range("C10").Value = 13680.9409932789
$ws.Range("B10").Formula = "=C10/C2"
$ws.Range("A11").Value = "BMsimVFI"

# --- Row 13 : ILAsolveLIN result now references column D instead of C ---
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = [double]"1.32821717998012E-2"
$ws.Range("B13").Formula = "=D13/D2"

# --- Row 14 / 15 : ILAsolveVFI <-> ILAsolveGSSA ---
# Row 14 becomes ILAsolveGSSA and now carries the result that used to sit on row 15
$ws.Range("A14").Value = "ILAsolveGSSA"
$ws.Range("F14").Value = 62.465973368613
$ws.Range("B14").Formula = "=F14/F2"
$ws.Range("B14").NumberFormat = "0.000"
# Row 15 becomes ILAsolveVFI, with its old B15/F15 numbers removed entirely
$ws.Range("A15").Value = "ILAsolveVFI"
$ws.Range("B15").Clear()
$ws.Range("F15").Clear()

# --- Row 16 : ILAsimLIN result now references column D instead of C ---
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = 10490.443206468901
$ws.Range("B16").Formula = "=D16/D2"

# --- Row 17 / 18 : ILAsimVFI <-> ILAsimGSSA ---
$ws.Range("A17").Value = "ILAsimGSSA"
$ws.Range("F17").Value = 12422.3521871921
$ws.Range("B17").Formula = "=F17/F2"
$ws.Range("A18").Value = "ILAsimVFI"
$ws.Range("B18").ClearContents()
$ws.Range("F18").Clear()

# --- Row 19 was a stray empty styled cell (B19) -- drop it completely ---
$ws.Range("B19").Clear()

# --- Row 20 : OLGsolveLIN result cleared back to no data (formula keeps referencing D20) ---
$ws.Range("D20").ClearContents()

# --- Row 21 / 22 : OLGsolveVFI <-> OLGsolveGSSA ---
$ws.Range("A21").Value = "OLGsolveGSSA"
$ws.Range("A22").Value = "OLGsolveVFI"

# --- Row 23 : OLGsimLIN result cleared back to no data (formula keeps referencing D23) ---
$ws.Range("D23").ClearContents()

# --- Row 24 / 25 : OLGsimVFI <-> OLGsimGSSA ---
$ws.Range("A24").Value = "OLGsimGSSA"
$ws.Range("A25").Value = "OLGsimVFI"
